$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.899.05'
$ws.Range("E2").Value = '  +1.39%  '

# Row 3
$ws.Range("D3").Value = '2.542.01'
$ws.Range("E3").Value = '  +0.87%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.32%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '315.98'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.63%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '96.60'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.08%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.573'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.17%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.26%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.537'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +1.48%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '35.68'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.17%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0810'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.83%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '7.48'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.33%  '

# Row 13
$ws.Range("E13").Value = '  -5.20%  '

# Row 14
$ws.Range("D14").Value = '2.916.88'
$ws.Range("E14").Value = '  +0.19%  '

# Row 15
$ws.Range("D15").Value = '2.501.46'
$ws.Range("E15").Value = '  -2.74%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '15.03'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.50%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.847'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.28%  '

# Row 18
$ws.Range("D18").Value = '42.739.84'
$ws.Range("E18").Value = '  +0.72%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.85'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +4.22%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '12.52'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.26%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0958'
$ws.Range("E21").Value = '  +0.12%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '69.49'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.36%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '253.25'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.42%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.95'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.81%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.05'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.75%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '26.56'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.10%  '

# Row 27
$ws.Range("E27").Value = '  +0.17%  '

# Row 28
$ws.Range("E28").Value = '  +1.23%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '40.41'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +4.34%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '10.24'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.68%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '5.83'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.16%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '156.32'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.23%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '19.57'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +4.40%  '

# Row 34
$ws.Range("E34").Value = '  +3.17%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0799'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +2.57%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.09'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.06%  '

# Row 37
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.28'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.56%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.111'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.97%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.43'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +4.10%  '

# Row 40
$ws.Range("E40").Value = '  +0.08%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '22.05'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -6.80%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '3.81'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.74%  '

# Row 43
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.0303'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.62%  '

# Row 44
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.50%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '3.24'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.70%  '

# Row 46
$ws.Range("D46").Value = '1.982.62'
$ws.Range("E46").Value = '  -1.63%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '84.40'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.36%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '8.99'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.24%  '

# Row 49
$ws.Range("D49").Value = '2.772.79'
$ws.Range("E49").Value = '  +0.28%  '

# Row 50
$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '74.34'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +2.91%  '

# Row 51
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '104.36'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +2.83%  '
